$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.5
$ws.Range("D2").Value = 0.42

$ws.Range("B3").Value = 12.92
$ws.Range("C3").Value = 12.5
$ws.Range("D3").Value = 8.75

$ws.Range("B4").Value = 21.67
$ws.Range("C4").Value = 22.5
$ws.Range("D4").Value = 19.58

$ws.Range("B5").Value = 42.08
$ws.Range("C5").Value = 40.83
$ws.Range("D5").Value = 39.58

$ws.Range("B6").Value = 56.25
$ws.Range("C6").Value = 52.5
$ws.Range("D6").Value = 51.67
